$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 181, shifting rows 181:183 down to 182:184
$ws.Rows.Item(181).Insert()

# Fill in the new row 181 data (copy of the "template" row, new date + new prices)
$ws.Cells.Item(181, 1).Value = 4
$ws.Cells.Item(181, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(181, 3).Value = "Los Lagos"
$ws.Cells.Item(181, 4).Value = 44448
$ws.Cells.Item(181, 5).Value = 10
$ws.Cells.Item(181, 6).Value = 100112023
$ws.Cells.Item(181, 7).Value = "Brócoli"
$ws.Cells.Item(181, 8).Value = "Sin especificar"
$ws.Cells.Item(181, 9).Value = "Primera"
$ws.Cells.Item(181, 10).Value = 500
$ws.Cells.Item(181, 11).Value = 1200
$ws.Cells.Item(181, 12).Value = 1200
$ws.Cells.Item(181, 13).Value = 1200
$ws.Cells.Item(181, 14).Value = "$/unidad"
$ws.Cells.Item(181, 15).Value = "Región Metropolitana"
$ws.Cells.Item(181, 16).Value = 1200
$ws.Cells.Item(181, 17).Value = 1
$ws.Cells.Item(181, 18).Value = "Hortaliza"
